$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53.. down by one.
$ws.Rows("53:53").Insert()

# Populate the newly inserted row 53 with the new record.
$ws.Range("A53").Value = 10
$ws.Range("B53").Value = "Vega Modelo de Temuco"
$ws.Range("C53").Value = "La Araucanía"
$ws.Range("D53").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D53").Value = 45033
$ws.Range("E53").Value = 9
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100107
$ws.Range("H53").Value = "Otros"
$ws.Range("I53").Value = 100107011
$ws.Range("J53").Value = "Tuna"
$ws.Range("K53").Value = "Sin especificar"
$ws.Range("L53").Value = "Primera"
$ws.Range("M53").Value = 150
$ws.Range("N53").Value = 22000
$ws.Range("O53").Value = 22000
$ws.Range("P53").Value = 22000
$ws.Range("Q53").Value = '$/caja 16 kilos'
$ws.Range("R53").Value = "Provincia de Los Andes"
$ws.Range("S53").Value = 1375
$ws.Range("T53").Value = 16
